# Gantt chart update: extend "Integration & Debugging", push back the
# "Extension: Simulator" / extension tasks by a week, rename the
# "Extension: M, A, F, D" task to include the base ISA, drop the two
# optional (IEEE754 / Extension Q&Help) rows, and shorten the final
# integration task.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Gantt Chart")

# 1) Row 29 "Integration & Debugging" - duration grows from 21 to 28 days
$ws.Range("E29").Value = 28

# 2) Row 30 "Extension: Simulator" - start date slips by one week
$ws.Range("C30").Value = 43177

# 3) Row 31 "Extension: M, A, F, D" -> rename (keep dates etc, just text)
#    and slip its start date by one week too
$ws.Range("B31").Value = "Base: RV32I, RV64I; Extension: M, A, F, D"
$ws.Range("C31").Value = 43177

# 4) Remove the two now-obsolete "(Optional) ..." rows (old rows 32 & 33)
$ws.Rows("32:33").Delete()

# 5) What used to be row 34 ("Final Integration & Debugging") is now row 32:
#    push its start date back a week and shrink its duration 21 -> 14 days
$ws.Range("C32").Value = 43191
$ws.Range("E32").Value = 14

$wb.Save()
